# re-calculate RQ1 up through google
# Adds a new recalculated/sorted breakdown (label + count) in columns H & K
# for rows 4-11, matching the original B-column categories in rows 4-11.
# Cells are populated in the same order the new unique label strings were
# first entered, so the generated shared-string table ordering matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entryOrder = @(
    @{ Row = 7;  Label = "incorrect program logic";     Count = 34 },
    @{ Row = 4;  Label = "algorithmic flakiness";       Count = 8 },
    @{ Row = 6;  Label = "environment";                 Count = 4 },
    @{ Row = 8;  Label = "async wait";                  Count = 20 },
    @{ Row = 9;  Label = "unordered collections";       Count = 2 },
    @{ Row = 10; Label = "concurrency";                 Count = 14 },
    @{ Row = 5;  Label = "memory";                      Count = 6 },
    @{ Row = 11; Label = "arithmetic/bit operations";   Count = 1 }
)

foreach ($item in $entryOrder) {
    $ws.Cells.Item($item.Row, 8).Value = $item.Label   # column H
    $ws.Cells.Item($item.Row, 11).Value = $item.Count  # column K
}

# Update the selected cell to match the post-edit state
$ws.Range("G11").Select()
